$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 216, shifting existing rows 216:248 down to 217:249
$ws.Rows("216:216").Insert()

# Populate the new row 216 with the new weekly price record
$ws.Range("A216").Value = 9
$ws.Range("B216").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C216").Value = "Metropolitana"
$ws.Range("D216").Value = 44776
$ws.Range("E216").Value = 13
$ws.Range("F216").Value = 100112026
$ws.Range("G216").Value = "Haba"
$ws.Range("H216").Value = "Sin especificar"
$ws.Range("I216").Value = "Primera"
$ws.Range("J216").Value = 52
$ws.Range("K216").Value = 18000
$ws.Range("L216").Value = 20000
$ws.Range("M216").Value = 19000
$ws.Range("N216").Value = "$/saco 25 kilos"
$ws.Range("O216").Value = "Región de Coquimbo"
$ws.Range("P216").Value = 760
$ws.Range("Q216").Value = 25
$ws.Range("R216").Value = "Hortaliza"
